$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.756.42'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').Value = '3.213.49'
$ws.Range('E3').Value = '  -3.05%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = '''576.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.49%  '
$ws.Range('D6').Value = '''141.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -12.77%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '3.213.12'
$ws.Range('E8').Value = '  -3.23%  '
$ws.Range('D9').Value = '''0.521'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -9.22%  '
$ws.Range('E10').Value = '  -11.66%  '
$ws.Range('D11').Value = '''6.43'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.39%  '
$ws.Range('D12').Value = '''0.478'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.61%  '
$ws.Range('D13').Value = '''0.0000231'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.83%  '
$ws.Range('D14').Value = '''36.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -12.16%  '
$ws.Range('D15').Value = '3.723.46'
$ws.Range('E15').Value = '  -4.12%  '
$ws.Range('D16').Value = '66.682.15'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').Value = '3.220.00'
$ws.Range('E17').Value = '  -4.60%  '
$ws.Range('E18').Value = '  -5.33%  '
$ws.Range('D19').Value = '''6.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.80%  '
$ws.Range('D20').Value = '''496.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -10.21%  '
$ws.Range('D21').Value = '''14.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -11.46%  '
$ws.Range('D22').Value = '''0.716'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.97%  '
$ws.Range('D23').Value = '''7.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -12.52%  '
$ws.Range('D24').Value = '''81.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.19%  '
$ws.Range('D25').Value = '''12.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.29%  '
$ws.Range('D26').Value = '''0.996'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').Value = '''3.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -11.58%  '
$ws.Range('E28').Value = '  -9.61%  '
$ws.Range('D29').Value = '''27.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.99%  '
$ws.Range('D30').Value = '''7.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.42%  '
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').Value = '''2.48'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('D34').Value = '''6.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -15.54%  '
$ws.Range('D35').Value = '''54.34'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').Value = '''5.30'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -13.34%  '
$ws.Range('D37').Value = '''485.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.74%  '
$ws.Range('D38').Value = '''0.0412'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.43%  '
$ws.Range('D39').Value = '''0.0809'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.06%  '
$ws.Range('E40').Value = '  -11.08%  '
$ws.Range('D41').Value = '''8.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -13.95%  '
$ws.Range('D42').Value = '2.824.73'
$ws.Range('E42').Value = '  -7.16%  '
$ws.Range('D43').Value = '''2.48'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.79%  '
$ws.Range('D44').Value = '''0.250'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.05%  '
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '''2.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.43%  '
$ws.Range('D47').Value = '''24.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -13.03%  '
$ws.Range('D48').Value = '''121.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.60%  '
$ws.Range('D49').Value = '0.0₃0522'
$ws.Range('E49').Value = '  -13.43%  '
$ws.Range('E50').Value = '  -8.97%  '
$ws.Range('E51').Value = '  -18.51%  '
